$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.905.65"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.299.49"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.510"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.118"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "2.658.60"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "2.295.91"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.784"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "42.842.24"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.72%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0692"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("D43").Value = "2.008.51"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "2.525.83"
$ws.Range("E51").Value = "  -0.20%  "
